# Fix Num of Cluster (100 data)
# Insert a new row at row 6 of sheet "_3_3" (the Section3 / Machinery_Part sheet).
# This pushes the former rows 6-10 down to rows 7-11, and the new row 6
# carries a new "I-beam with two (2) sets of trolleys / △2." entry whose
# data columns are mostly "0" (except column BF, which stays "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_3_3")

# Insert a new blank row above the current row 6; existing rows 6-10 shift to 7-11.
$ws.Rows("6").Insert()

# Column A: descriptive text for the new row.
$ws.Cells.Item(6, 1).Value = "I-beam with two (2) sets of trolleys`n△2."

# Columns B (2) through CX (102): default to "0" (stored as text, like the
# rest of the sheet, not as a number). Column BF (58) keeps the value "1"
# instead of "0".
for ($col = 2; $col -le 102; $col++) {
    if ($col -eq 58) {
        $textValue = "1"
    } else {
        $textValue = "0"
    }
    $cell = $ws.Cells.Item(6, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $textValue
    $cell.Style = "Normal"
}
